$d = $word.ActiveDocument

function Replace-Span([string]$findText, [string]$newText) {
    # Locate the exact span, then force a genuine content rewrite (even when
    # $newText happens to equal the old text) by routing through a unique
    # placeholder.  This also coalesces every run / proofErr marker inside
    # the span into a single clean run, mirroring what Word does when a
    # user retypes across an existing selection.
    $t = $d.Content.Text
    $idx = $t.IndexOf($findText)
    if ($idx -lt 0) {
        throw "Replace-Span: text not found: $findText"
    }
    $r = $d.Range($idx, $idx + $findText.Length)
    $marker = "@@MARKER_" + [guid]::NewGuid().ToString("N") + "@@"
    $r.Text = $marker
    $t2 = $d.Content.Text
    $idx2 = $t2.IndexOf($marker)
    $r2 = $d.Range($idx2, $idx2 + $marker.Length)
    $r2.Text = $newText
}

# ---------------------------------------------------------------------
# 1) "Units" paragraph: celsius -> Celsius, and grammar fix around
#    "will use take ... if necessary convert".
# ---------------------------------------------------------------------
Replace-Span `
    "supplied in degrees celsius.  There exist functions for converting between units, and many of the output or input functions have versions which take a UNITS input.  The functions with this input will use take the values specified and if necessary convert them to the appropriate units.  " `
    "supplied in degrees Celsius.  There exist functions for converting between units, and many of the output or input functions have versions which take a UNITS input.  The functions with this input will take the values specified and, if necessary, convert them to the appropriate units.  "

# ---------------------------------------------------------------------
# 2) "Interior space temperature depression" paragraph: strip the stray
#    proofErr spell/gram markers around "Ecotope's" and the "Because
#    HPWH's ..." sentence (text itself is unchanged).
# ---------------------------------------------------------------------
Replace-Span `
    " was originally written for SEEM, Ecotope's single-zone building energy use simulation engine.  Because HPWH's remove energy from the air, there is the possibility that they will depress the temperature of their local environment, thus decreasing their performance.  This effect was measured in a field study, with the average temperature depression being approximately 4.5 F with a " `
    " was originally written for SEEM, Ecotope's single-zone building energy use simulation engine.  Because HPWH's remove energy from the air, there is the possibility that they will depress the temperature of their local environment, thus decreasing their performance.  This effect was measured in a field study, with the average temperature depression being approximately 4.5 F with a "

# ---------------------------------------------------------------------
# 3) "Hysteresis" paragraph: celsius -> Celsius (and drop the stray
#    spellStart/gramStart proofErr markers around it).
# ---------------------------------------------------------------------
Replace-Span `
    ", are differential celsius degrees" `
    ", are differential Celsius degrees"

# ---------------------------------------------------------------------
# 4) Fatal-errors bullet list: drop the "_GoBack" bookmark from in front
#    of "A value other than 1 for minutesPerStep ...", fix the "is the
#    verbosity" typo to "if the verbosity", and re-plant the bookmark at
#    its new spot right after "message if".
# ---------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

Replace-Span `
    "This kind of fatal error will write out an informative error message is the verbosity is not set to" `
    "This kind of fatal error will write out an informative error message if the verbosity is not set to"

$t = $d.Content.Text
$anchor = "This kind of fatal error will write out an informative error message if"
$idx = $t.IndexOf($anchor)
$insertPos = $idx + $anchor.Length
$r = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $r)
